{"js": "// Update the division problems in the practice-sheet table.\n// Each <w:t> run holding a \"NNN\u00f7N=\" expression is replaced in place so\n// that the existing run formatting (font/size) is preserved.\nconst replacements = [\n  [\"681\u00f77=\", \"326\u00f79=\"],\n  [\"435\u00f72=\", \"355\u00f77=\"],\n  [\"154\u00f74=\", \"304\u00f78=\"],\n  [\"895\u00f76=\", \"366\u00f74=\"],\n  [\"565\u00f78=\", \"567\u00f73=\"],\n  [\"714\u00f72=\", \"809\u00f73=\"],\n  [\"720\u00f79=\", \"468\u00f76=\"],\n  [\"279\u00f77=\", \"855\u00f73=\"],\n  [\"222\u00f74=\", \"955\u00f76=\"],\n  [\"531\u00f77=\", \"880\u00f73=\"],\n  [\"980\u00f77=\", \"203\u00f74=\"],\n  [\"701\u00f72=\", \"685\u00f77=\"],\n  [\"349\u00f76=\", \"865\u00f78=\"],\n  [\"145\u00f77=\", \"693\u00f72=\"],\n  [\"508\u00f74=\", \"315\u00f72=\"],\n  [\"723\u00f77=\", \"660\u00f72=\"],\n  [\"582\u00f75=\", \"764\u00f73=\"],\n  [\"286\u00f78=\", \"798\u00f77=\"],\n  [\"332\u00f79=\", \"880\u00f73=\"],\n  [\"812\u00f79=\", \"114\u00f79=\"],\n  [\"773\u00f79=\", \"367\u00f75=\"],\n  [\"603\u00f78=\", \"495\u00f72=\"],\n  [\"293\u00f73=\", \"938\u00f74=\"],\n  [\"751\u00f79=\", \"342\u00f73=\"],\n  [\"409\u00f77=\", \"405\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the division problems in the practice-sheet table.\n# Each Find/Replace targets one \"NNN\u00f7N=\" run; wildcards are off and the\n# match is whole-text, so only the intended run is touched and its\n# formatting (font/size) is left untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"681\u00f77=\", \"326\u00f79=\"),\n    @(\"435\u00f72=\", \"355\u00f77=\"),\n    @(\"154\u00f74=\", \"304\u00f78=\"),\n    @(\"895\u00f76=\", \"366\u00f74=\"),\n    @(\"565\u00f78=\", \"567\u00f73=\"),\n    @(\"714\u00f72=\", \"809\u00f73=\"),\n    @(\"720\u00f79=\", \"468\u00f76=\"),\n    @(\"279\u00f77=\", \"855\u00f73=\"),\n    @(\"222\u00f74=\", \"955\u00f76=\"),\n    @(\"531\u00f77=\", \"880\u00f73=\"),\n    @(\"980\u00f77=\", \"203\u00f74=\"),\n    @(\"701\u00f72=\", \"685\u00f77=\"),\n    @(\"349\u00f76=\", \"865\u00f78=\"),\n    @(\"145\u00f77=\", \"693\u00f72=\"),\n    @(\"508\u00f74=\", \"315\u00f72=\"),\n    @(\"723\u00f77=\", \"660\u00f72=\"),\n    @(\"582\u00f75=\", \"764\u00f73=\"),\n    @(\"286\u00f78=\", \"798\u00f77=\"),\n    @(\"332\u00f79=\", \"880\u00f73=\"),\n    @(\"812\u00f79=\", \"114\u00f79=\"),\n    @(\"773\u00f79=\", \"367\u00f75=\"),\n    @(\"603\u00f78=\", \"495\u00f72=\"),\n    @(\"293\u00f73=\", \"938\u00f74=\"),\n    @(\"751\u00f79=\", \"342\u00f73=\"),\n    @(\"409\u00f77=\", \"405\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
